# New crime data collected - weekly CompStat update (17th Precinct)
# Updates the report header (volume/date range) and the weekly crime-stat
# table (rows 15-27) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 30   Number  22" -> "...  23"
# and the reporting week date range.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# ---------------------------------------------------------------------
# Some table cells hold the literal text "0" or "***.*" (the report
# template's way of saying "no data" / "undefined % change") instead of
# a real number, while sibling cells hold genuine numbers. This week's
# data flips several cells each way, so both directions are needed:
#
#  * numeric -> text: typing a digit-only string like "0" auto-converts
#    to a number (same as real Excel) unless entered with a leading
#    apostrophe (quote-prefix) to force text. Non-numeric-looking text
#    such as "***.*" does not need the apostrophe. Either way, the
#    quote-prefix / previous numeric format can leave the cell's format
#    not matching its text sibling cells, so we copy the number format
#    from an already-correct "no data" cell in row 14 (untouched by
#    this week's edits) to line it back up.
#  * text -> numeric: after writing the number we restore the normal
#    numeric display format (plain count, or signed one-decimal percent)
#    that the rest of the column already uses.
# ---------------------------------------------------------------------

$countFmt = "#,##0"
$pctFmt = "#,##0.0;""-""#,##0.0"

# Row 15 - Rape
$ws.Range("G15").Value = "'0"
$ws.Range("H15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4122)

# Row 16 - Robbery
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("D16").NumberFormat = $countFmt
$ws.Range("E16").Value = 0
$ws.Range("E16").NumberFormat = $pctFmt
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 35
$ws.Range("L16").Value = 34.615384615384
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = -86.842105263157

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 2
$ws.Range("D17").NumberFormat = $countFmt
$ws.Range("E17").Value = -50
$ws.Range("E17").NumberFormat = $pctFmt
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 44
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = -6.382978723404
$ws.Range("L17").Value = -15.384615384615
$ws.Range("M17").Value = 37.5
$ws.Range("N17").Value = -43.589743589743

# Row 18 - Burglary
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = -37.362637362637
$ws.Range("L18").Value = 3.636363636363
$ws.Range("M18").Value = 7.547169811320
$ws.Range("N18").Value = -88.271604938271

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 46.666666666666
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = 14.285714285714
$ws.Range("I19").Value = 333
$ws.Range("J19").Value = 284
$ws.Range("K19").Value = 17.253521126760
$ws.Range("L19").Value = 56.338028169014
$ws.Range("M19").Value = 5.379746835443
$ws.Range("N19").Value = -64.574468085106

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
$ws.Range("C20").NumberFormat = $countFmt
$ws.Range("E20").Value = 100
$ws.Range("I20").Value = 24
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = 41.176470588235
$ws.Range("L20").Value = -4
$ws.Range("M20").Value = 118.181818181818
$ws.Range("N20").Value = -91.808873720136

# Row 21 - TOTAL
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 7.407407407407
$ws.Range("F21").Value = 88
$ws.Range("H21").Value = -4.347826086956
$ws.Range("I21").Value = 495
$ws.Range("J21").Value = 478
$ws.Range("K21").Value = 3.556485355648
$ws.Range("L21").Value = 32.707774798927
$ws.Range("M21").Value = 11.738148984198
$ws.Range("N21").Value = -76.063829787234

# Row 22 - Transit
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = $countFmt
$ws.Range("D22").Value = 2
$ws.Range("D22").NumberFormat = $countFmt
$ws.Range("E22").Value = -50
$ws.Range("E22").NumberFormat = $pctFmt
$ws.Range("G22").Value = 2
$ws.Range("G22").NumberFormat = $countFmt
$ws.Range("H22").Value = 50
$ws.Range("H22").NumberFormat = $pctFmt
$ws.Range("I22").Value = 14
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 40
$ws.Range("L22").Value = 7.692307692307
$ws.Range("M22").Value = -17.647058823529

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = -19.354838709677
$ws.Range("I24").Value = 417
$ws.Range("J24").Value = 532
$ws.Range("K24").Value = -21.616541353383
$ws.Range("L24").Value = -6.292134831460
$ws.Range("M24").Value = 47.349823321554

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 20
$ws.Range("H25").Value = 17.647058823529
$ws.Range("I25").Value = 122
$ws.Range("J25").Value = 107
$ws.Range("K25").Value = 14.018691588785
$ws.Range("L25").Value = 79.411764705882
$ws.Range("M25").Value = 16.190476190476

# Row 26 - UCR Rape*
$ws.Range("G26").Value = "'0"
$ws.Range("H26").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H26").PasteSpecial(-4122)

# Row 27 - Other Sex Crimes
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = $countFmt
$ws.Range("E27").Value = 100
$ws.Range("E27").NumberFormat = $pctFmt
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = 27.272727272727
$ws.Range("L27").Value = -9.677419354838
